$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 369.69232
$ws.Range("I9").Value = 189.33333
$ws.Range("K9").Value = 189.33333
$ws.Range("M9").Value = -20.33332999999999
$ws.Range("H15").Value = 2422.7307
$ws.Range("I15").Value = 2422.7307
$ws.Range("K15").Value = 7268.1921
$ws.Range("M15").Value = -7099.1921
$ws.Range("H17").Value = 808.6667
$ws.Range("J17").Value = 808.6667
$ws.Range("L17").Value = 2426.0001
$ws.Range("N17").Value = -2762.0001
$ws.Range("H39").Value = 1153.1154
$ws.Range("I39").Value = 680.2273
$ws.Range("J39").Value = 3754
$ws.Range("K39").Value = 2040.6819
$ws.Range("L39").Value = 11262
$ws.Range("M39").Value = -1744.6819
$ws.Range("N39").Value = -11854
$ws.Range("H117").Value = 60000
$ws.Range("J117").Value = 60000
$ws.Range("L117").Value = 60000
$ws.Range("N117").Value = -69178
$ws.Range("H132").Value = 5090.381
$ws.Range("I132").Value = 3944.9211
$ws.Range("J132").Value = 15972.25
$ws.Range("K132").Value = 11834.7633
$ws.Range("L132").Value = 47916.75
$ws.Range("M132").Value = -9304.763300000001
$ws.Range("N132").Value = -52976.75
$ws.Range("H137").Value = 4361.375
$ws.Range("I137").Value = 4413.143
$ws.Range("J137").Value = 3999
$ws.Range("K137").Value = 13239.429
$ws.Range("L137").Value = 11997
$ws.Range("M137").Value = -10689.429
$ws.Range("N137").Value = -17097

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 5999
$ws.Range("I102").Value = 5999
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 5999
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -4377
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 1594.35
$ws.Range("I122").Value = 1478.2632
$ws.Range("K122").Value = 4434.7896
$ws.Range("M122").Value = -1984.7896
$ws.Range("H131").Value = 79987
$ws.Range("J131").Value = 79987
$ws.Range("L131").Value = 79987
$ws.Range("N131").Value = -90067
$ws.Range("H132").Value = 33182.312
$ws.Range("I132").Value = 69600.53
$ws.Range("J132").Value = 5868.65
$ws.Range("K132").Value = 208801.59
$ws.Range("L132").Value = 17605.95
$ws.Range("M132").Value = -206271.59
$ws.Range("N132").Value = -22665.95

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 49999.5
$ws.Range("J15").Value = 49999.5
$ws.Range("L15").Value = 49999.5
$ws.Range("N15").Value = -50453.5
$ws.Range("H107").Value = 3611.25
$ws.Range("I107").Value = 2669
$ws.Range("J107").Value = 3925.3333
$ws.Range("K107").Value = 2669
$ws.Range("L107").Value = 3925.3333
$ws.Range("M107").Value = -749
$ws.Range("N107").Value = -7765.3333
$ws.Range("H118").Value = 79997.5
$ws.Range("J118").Value = 79997.5
$ws.Range("L118").Value = 79997.5
$ws.Range("N118").Value = -83311.5
$ws.Range("H127").Value = 74999
$ws.Range("J127").Value = 74999
$ws.Range("L127").Value = 74999
$ws.Range("N127").Value = -84919
$ws.Range("H133").Value = 99497.5
$ws.Range("J133").Value = 99497.5
$ws.Range("L133").Value = 99497.5
$ws.Range("N133").Value = -109617.5
$ws.Range("H135").Value = 76815.57000000001
$ws.Range("J135").Value = 76815.57000000001
$ws.Range("L135").Value = 76815.57000000001
$ws.Range("N135").Value = -86955.57000000001
$ws.Range("H137").Value = 69298
$ws.Range("J137").Value = 69298
$ws.Range("L137").Value = 69298
$ws.Range("N137").Value = -79498

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 977.381
$ws.Range("I122").Value = 940.13336
$ws.Range("K122").Value = 2820.40008
$ws.Range("M122").Value = -370.4000800000003
$ws.Range("H132").Value = 3108
$ws.Range("I132").Value = 2912
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 8736
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -6206
$ws.Range("N132").Value = -15560

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 599.875
$ws.Range("J33").Value = 874.75
$ws.Range("L33").Value = 5248.5
$ws.Range("N33").Value = -5814.5
$ws.Range("H95").Value = 6999.8335
$ws.Range("J95").Value = 6999.8335
$ws.Range("L95").Value = 20999.5005
$ws.Range("N95").Value = -25117.5005
$ws.Range("H105").Value = 8967.813
$ws.Range("J105").Value = 8967.813
$ws.Range("L105").Value = 26903.439
$ws.Range("N105").Value = -32145.439
$ws.Range("H107").Value = 400874
$ws.Range("J107").Value = 435658.56
$ws.Range("L107").Value = 1306975.68
$ws.Range("N107").Value = -1310815.68
$ws.Range("H114").Value = 754.4211
$ws.Range("I114").Value = 656
$ws.Range("J114").Value = 923.1429000000001
$ws.Range("K114").Value = 1968
$ws.Range("L114").Value = 2769.4287
$ws.Range("M114").Value = 1286
$ws.Range("N114").Value = -9277.4287
$ws.Range("H132").Value = 3078.111
$ws.Range("I132").Value = 3524.5
$ws.Range("K132").Value = 31720.5
$ws.Range("M132").Value = -29190.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 6625000
$ws.Range("J7").Value = 6625000
$ws.Range("L7").Value = 6625000
$ws.Range("N7").Value = -6625224
$ws.Range("H8").Value = 6625000
$ws.Range("J8").Value = 6625000
$ws.Range("L8").Value = 6625000
$ws.Range("N8").Value = -6625278
$ws.Range("H94").Value = 69999
$ws.Range("J94").Value = 69999
$ws.Range("L94").Value = 69999
$ws.Range("N94").Value = -71351
$ws.Range("H102").Value = 1931.5834
$ws.Range("I102").Value = 1931.5834
$ws.Range("K102").Value = 1931.5834
$ws.Range("M102").Value = -309.5834
$ws.Range("H122").Value = 3062.0908
$ws.Range("I122").Value = 2381.5557
$ws.Range("J122").Value = 6124.5
$ws.Range("K122").Value = 7144.6671
$ws.Range("L122").Value = 18373.5
$ws.Range("M122").Value = -4694.6671
$ws.Range("N122").Value = -23273.5
$ws.Range("H132").Value = 75331.28999999999
$ws.Range("I132").Value = 87553.25
$ws.Range("J132").Value = 1999.5
$ws.Range("K132").Value = 262659.75
$ws.Range("L132").Value = 5998.5
$ws.Range("M132").Value = -260129.75
$ws.Range("N132").Value = -11058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 3010
$ws.Range("I35").Value = 3010
$ws.Range("K35").Value = 3010
$ws.Range("M35").Value = -2674
$ws.Range("H122").Value = 3502
$ws.Range("I122").Value = 3502
$ws.Range("K122").Value = 10506
$ws.Range("M122").Value = -8056
$ws.Range("H132").Value = 54204.957
$ws.Range("I132").Value = 58846.24
$ws.Range("K132").Value = 176538.72
$ws.Range("M132").Value = -174008.72
$ws.Range("H136").Value = 9278
$ws.Range("I136").Value = 8762.286
$ws.Range("K136").Value = 26286.858
$ws.Range("M136").Value = -23736.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3405
$ws.Range("I81").Value = 3405
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 6810
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -5749
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 3405
$ws.Range("I84").Value = 3405
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 34050
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -28746
$ws.Range("N84").ClearContents()
$ws.Range("H113").Value = 4029
$ws.Range("I113").Value = 1952.1428
$ws.Range("K113").Value = 5856.428400000001
$ws.Range("M113").Value = -3686.428400000001
$ws.Range("H122").Value = 7053.4546
$ws.Range("I122").Value = 7949.421
$ws.Range("J122").Value = 1379
$ws.Range("K122").Value = 23848.263
$ws.Range("L122").Value = 4137
$ws.Range("M122").Value = -21398.263
$ws.Range("N122").Value = -9037
$ws.Range("H136").Value = 2780.0386
$ws.Range("I136").Value = 1767.409
$ws.Range("K136").Value = 5302.227000000001
$ws.Range("M136").Value = -2752.227000000001
